$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("coor")

# New hospital rows 13-16 referencing the lookup table in A1:E9, plus a totals row 17.
$ws.Range("B13").Value = "h6"
$ws.Range("B14").Value = "h3"
$ws.Range("B15").Value = "h2"
$ws.Range("B16").Value = "h7"

$ws.Range("C13").Formula = '=VLOOKUP(B13,$A$1:$E$9,4,)'
$ws.Range("D13").Formula = '=VLOOKUP(B13,$A$1:$E$9,5,)'

$ws.Range("C14:C16").Formula = '=VLOOKUP(B14,$A$1:$E$9,4,)'
$ws.Range("D14:D16").Formula = '=VLOOKUP(B14,$A$1:$E$9,5,)'

$ws.Range("C17").Formula = '=SUM(C13:C16)'
$ws.Range("D17").Formula = '=SUM(D13:D16)'

# Make "coor" the active sheet/tab and select D14 (matches the new selection in the file).
$ws.Activate() | Out-Null
$ws.Select() | Out-Null
$ws.Range("D14").Select() | Out-Null
